$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 ("Tháng 7") data refresh: Notion export re-synced with updated
# numbers, bumping last_edited_time for this page (and a few sibling pages
# that happened to share the previous last_edited_time timestamp).

# Updated figures for Tháng 7 (row 13)
$ws.Range("W13").Value = 105486000
$ws.Range("AA13").Value = 216944000
$ws.Range("AE13").Value = 322430000
$ws.Range("AH13").Value = 269730000
$ws.Range("AK13").Value = 47
$ws.Range("AN13").Value = 52700000
$ws.Range("AQ13").Value = 305530000

# last_edited_time bump for the edited row plus the sibling rows that
# previously shared the same timestamp value.
$ws.Range("D5").Value = "2024-07-31T18:25:00.000Z"
$ws.Range("D8").Value = "2024-07-31T18:25:00.000Z"
$ws.Range("D13").Value = "2024-07-31T18:25:00.000Z"

$ws.Range("D4").Value = "2024-07-31T18:24:00.000Z"
$ws.Range("D6").Value = "2024-07-31T18:24:00.000Z"
$ws.Range("D12").Value = "2024-07-31T18:24:00.000Z"
